$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target values for rows 23-32, columns D, L, M, N, O, P, S
# (Rows shift down by one, row 23 receives a new record, row 32 receives
# the previously-missing record; all other columns are unchanged.)

$ws.Range("D23").Value = 44586
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 150
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 3000
$ws.Range("S23").Value = 1500

$ws.Range("D24").Value = 44524
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 180
$ws.Range("N24").Value = 4000
$ws.Range("O24").Value = 4000
$ws.Range("P24").Value = 4000
$ws.Range("S24").Value = 2000

$ws.Range("D25").Value = 44530
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 4000
$ws.Range("O25").Value = 4000
$ws.Range("P25").Value = 4000
$ws.Range("S25").Value = 2000

$ws.Range("D26").Value = 44530
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 150
$ws.Range("N26").Value = 3600
$ws.Range("O26").Value = 3600
$ws.Range("P26").Value = 3600
$ws.Range("S26").Value = 1800

$ws.Range("D27").Value = 44582
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 3200
$ws.Range("O27").Value = 3200
$ws.Range("P27").Value = 3200
$ws.Range("S27").Value = 1600

$ws.Range("D28").Value = 44235
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 3000
$ws.Range("O28").Value = 3000
$ws.Range("P28").Value = 3000
$ws.Range("S28").Value = 1500

$ws.Range("D29").Value = 44516
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 5000
$ws.Range("O29").Value = 5000
$ws.Range("P29").Value = 5000
$ws.Range("S29").Value = 2500

$ws.Range("D30").Value = 44552
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 180
$ws.Range("N30").Value = 4000
$ws.Range("O30").Value = 4000
$ws.Range("P30").Value = 4000
$ws.Range("S30").Value = 2000

$ws.Range("D31").Value = 44211
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 40
$ws.Range("N31").Value = 2800
$ws.Range("O31").Value = 2800
$ws.Range("P31").Value = 2800
$ws.Range("S31").Value = 1400

$ws.Range("D32").Value = 44211
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 30
$ws.Range("N32").Value = 2600
$ws.Range("O32").Value = 2600
$ws.Range("P32").Value = 2600
$ws.Range("S32").Value = 1300
